$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("U7").Value = 1.7
$ws.Range("V7").Value = 2.05
$ws.Range("AE7").Value = 15

$ws.Range("M8").Value = 1.03
$ws.Range("O8").Value = 1.25

$ws.Range("M9").Value = 1.03
$ws.Range("O9").Value = 1.22

$ws.Range("U15").Value = 1.47

$ws.Range("U16").Value = 1.5
$ws.Range("V16").Value = 2.37

$ws.Range("G17").Value = 1.53
$ws.Range("H17").Value = 4.2
$ws.Range("I17").Value = 5.75
$ws.Range("L17").Value = 6
$ws.Range("M17").Value = 1.05
$ws.Range("N17").Value = 11
$ws.Range("U17").Value = 1.87
$ws.Range("V17").Value = 1.77
$ws.Range("W17").Value = 6.5
$ws.Range("AH17").Value = 15
$ws.Range("AK17").Value = 67
$ws.Range("AW17").Value = 7.5

$ws.Range("M18").Value = 1.01
$ws.Range("O18").Value = 1.1

$ws.Range("M19").Value = 1.03
$ws.Range("O19").Value = 1.17

$ws.Range("M20").Value = 1.03
$ws.Range("O20").Value = 1.19

$ws.Range("M21").Value = 1.01
$ws.Range("O21").Value = 1.11

$ws.Range("O22").Value = 1.07

$ws.Range("U23").Value = 1.87
$ws.Range("V23").Value = 1.77

$ws.Range("U25").Value = 1.77
$ws.Range("V25").Value = 1.92

$ws.Range("U26").Value = 1.58

$ws.Range("V27").Value = 1.69

$ws.Range("U28").Value = 1.69

$ws.Range("U30").Value = 1.77
$ws.Range("V30").Value = 1.87

$ws.Range("U31").Value = 1.47

$ws.Range("U32").Value = 1.92
$ws.Range("V32").Value = 1.77

$ws.Range("K33").Value = 2.37
$ws.Range("U33").Value = 1.87
$ws.Range("V33").Value = 1.87

$ws.Range("U34").Value = 1.87
$ws.Range("V34").Value = 1.77

$ws.Range("J35").Value = 2.87
